$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos price/volume refresh (scraped values). Price column cells that
# look numeric must stay literal text (to preserve exact digits/trailing
# zeros as scraped), so force Text format before writing those.

$ws.Range("D2").Value = "29.114.69"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.835.43"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.34"
$ws.Range("E5").Value = "  +1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6342"
$ws.Range("E6").Value = "  +1.87%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07537"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2948"
$ws.Range("E9").Value = "  +1.44%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.05"
$ws.Range("E10").Value = "  +1.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07715"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D12").Value = "1.830.62"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.007"
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6725"
$ws.Range("E14").Value = "  +1.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "83.32"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009615"
$ws.Range("E16").Value = "  +5.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.093"
$ws.Range("E17").Value = "  +1.91%  "
$ws.Range("D18").Value = "29.130.13"
$ws.Range("E18").Value = "  +0.70%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.62"
$ws.Range("E19").Value = "  +2.46%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "227.33"
$ws.Range("E20").Value = "  +1.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.187"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.07%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "160.44"
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1418"
$ws.Range("E25").Value = "  +4.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.561"
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.97"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.502"
$ws.Range("E28").Value = "  +0.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.156"
$ws.Range("E29").Value = "  +2.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.077"
$ws.Range("E30").Value = "  +1.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.201"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05420"
$ws.Range("E32").Value = "  +4.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.862"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7471"
$ws.Range("E34").Value = "  +2.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.142"
$ws.Range("E35").Value = "  -0.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.658"
$ws.Range("E36").Value = "  +1.73%  "
$ws.Range("D37").Value = "1.248.52"
$ws.Range("E37").Value = "  -2.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01794"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.758"
$ws.Range("E39").Value = "  -0.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.671"
$ws.Range("E40").Value = "  +4.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9051"
$ws.Range("E41").Value = "  +1.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.002"
$ws.Range("E42").Value = "  +0.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "101.71"
$ws.Range("E43").Value = "  +0.39%  "
$ws.Range("D44").Value = "1.981.79"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000124"
$ws.Range("E45").Value = "  +4.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.35"
$ws.Range("E46").Value = "  +3.02%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5120"
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4075"
$ws.Range("E48").Value = "  +2.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.946"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05784"
$ws.Range("E50").Value = "  +0.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.774"
$ws.Range("E51").Value = "  +1.41%  "
